$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row for "Living With a SEAL" (Jesse Itzler)
$ws.Range("A8").Value = "Living With a SEAL"
$ws.Range("B8").Value = "Jesse Itzler"
# Copy the date formatting (style) from the row above, then set the raw
# serial values, so the new date cells reuse the existing style (no new
# number-format entries get added to styles.xml)
$ws.Range("C7:D7").Copy()
$ws.Range("C8:D8").PasteSpecial(-4122)
$ws.Range("C8").Value = 43840
$ws.Range("D8").Value = 43842

$ws.Range("E8").Value = "exercise;motivation;self-improvement"
$ws.Range("F8").Value = "Audio"
$ws.Range("G8").Value = "5 Hrs 19 Mins"

# Keep the "next empty row" selection convention used in this sheet
$ws.Range("A9").Select()
